$d = $word.ActiveDocument

# --- 1. "age" bookmark: shrink so it only covers "17" (not " YEARS OLD ") ---
$ageBm = $d.Bookmarks("age")
$ageStart = $ageBm.Start
$ageEnd = $ageBm.End
$findRange = $d.Range($ageStart, $ageEnd)
[void]$findRange.Find.Execute("17", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("age", $d.Range($ageStart, $findRange.End))

# --- 2. "purpose" bookmark: delete its contained text, then delete the bookmark itself ---
$purposeBm = $d.Bookmarks("purpose")
$purposeBm.Range.Text = ""
$d.Bookmarks("purpose").Delete()

# --- 3. "remark" bookmark: delete its contained text, but keep the (now empty) bookmark ---
$remarkBm = $d.Bookmarks("remark")
$remarkBm.Range.Text = ""

# "dateIssued" bookmark keeps its name/content; its w:id simply renumbers
# automatically as a consequence of removing the "purpose" bookmark above.
